# "export pdf control version"
#
# The paragraph that starts with "ตาราง … เป้าหมายทีม" (a placeholder table
# caption) is updated to read "ตารางที่ 1 เป้าหมายทีม" (table caption with an
# actual number filled in). The original paragraph is made up of 4 runs:
#   R1 "ตาราง"        (TH Sarabun New, cs, hint=cs)
#   R2 " "            (TH Sarabun New, cs, hint=cs)   <- identical rPr to R1
#   R3 "… "           (TH Sarabun New, no cs/hint)
#   R4 "เป้าหมายทีม"   (TH Sarabun New, cs, hint=cs)
# and needs to become 5 runs:
#   R1 "ตารา"         (unchanged rPr)
#   R2 "งที่ "         (unchanged rPr, identical to R1's)
#   R3 "1"            (unchanged rPr)
#   R4(new) " "       (same rPr as R3)
#   R5 "เป้าหมายทีม"   (unchanged, = old R4)
#
# iron_native's Range.Text setter merges the edited run together with any
# *directly adjacent* run that shares identical run-properties (this mirrors
# how Word itself treats a run of consecutive same-format runs as one
# editable span). Since R1/R2 share formatting (and, after the edit, the two
# pieces carved out of R3 must also share formatting), editing naively would
# collapse them back into a single <w:r>. To keep the runs distinct we
# briefly toggle a formatting flag (Bold) on the neighbour before touching
# text - that one-character difference blocks the auto-merge - and flip it
# back once both sides have their final text. Property-only writes never
# trigger the merge, so the final restore is safe.

$d = $word.ActiveDocument

# Locate the target paragraph (the placeholder table caption) defensively,
# rather than assuming a fixed paragraph index.
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.StartsWith("ตาราง")) {
        $target = $para.Range
        break
    }
}

$pStart = $target.Start

# --- R1 "ตาราง" / R2 " " --------------------------------------------------
# R2 is the single space right after "ตาราง" (5 characters in).
$r2 = $d.Range($pStart + 5, $pStart + 6)
if ($r2.Text -ne " ") { throw "unexpected run2 text: [$($r2.Text)]" }
$r2.Font.Bold = $true                      # differ from R1 so edits don't merge

$r1 = $d.Range($pStart, $pStart + 5)
if ($r1.Text -ne "ตาราง") { throw "unexpected run1 text" }
$r1.Text = "ตารา"                          # R1: "ตาราง" -> "ตารา"

$r2 = $d.Range($pStart + 4, $pStart + 5)
if ($r2.Text -ne " ") { throw "unexpected run2 text (2): [$($r2.Text)]" }
$r2.Text = "งที่ "                          # R2: " " -> "งที่ "

$r2 = $d.Range($pStart + 4, $pStart + 9)
$r2.Font.Bold = $false                     # restore formatting (no merge on property writes)

# --- R3 "… " -> "1" + new run " " -----------------------------------------
# Layout so far from $pStart: "ตารา"(4) + "งที่ "(5) + "… "(2) + ...
$ellipsis = $pStart + 9
$trail = $ellipsis + 1

$r3trail = $d.Range($trail, $trail + 1)
if ($r3trail.Text -ne " ") { throw "unexpected run3 trailing text: [$($r3trail.Text)]" }
$r3trail.Font.Bold = $true                 # differ from the "…" char so the split sticks

$r3main = $d.Range($ellipsis, $ellipsis + 1)
if ($r3main.Text -ne [char]0x2026) { throw "unexpected run3 text" }
$r3main.Text = "1"                         # R3: "…" -> "1"

$r3trail = $d.Range($ellipsis + 1, $ellipsis + 2)
$r3trail.Font.Bold = $false                # restore formatting -> stays its own run

$final = $target.Text
if (-not $final.StartsWith("ตารางที่ 1 เป้าหมายทีม")) { throw "final text mismatch: [$final]" }
Write-Output ("final paragraph text: " + $final)
